# Explicit wait and other modifications.
#
# Update the "Bills" reference number and the "Login" contact e-mail, then
# make "Login" the active/selected sheet (it was "Bills" before).

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")
$billsSheet = $wb.Worksheets.Item("Bills")

# Update the Bills reference number first (B2) so it keeps its existing
# shared-string slot; force it to stay text (leading apostrophe) so it
# keeps the quote-prefixed text style instead of becoming a number.
$billsSheet.Range("B2").Value = "'99974010169"

# Update the Login e-mail / id in A1.
$loginSheet.Range("A1").Value = "Cantt@mc.com"

# Make "Login" the active sheet/tab (previously "Bills" was active).
$loginSheet.Activate()
